# Insert a new data row at row 630 (pushing existing rows 630-735 down to
# 631-736), then populate the newly inserted row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(630).Insert()

$ws.Cells.Item(630, 1).Value = 5
$ws.Cells.Item(630, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(630, 3).Value = "Maule"
$ws.Cells.Item(630, 4).Value = 45180
$ws.Cells.Item(630, 5).Value = 7
$ws.Cells.Item(630, 6).Value = 100112043
$ws.Cells.Item(630, 7).Value = "Pepino ensalada"
$ws.Cells.Item(630, 8).Value = "Sin especificar"
$ws.Cells.Item(630, 9).Value = "Primera"
$ws.Cells.Item(630, 10).Value = 300
$ws.Cells.Item(630, 11).Value = 15000
$ws.Cells.Item(630, 12).Value = 15000
$ws.Cells.Item(630, 13).Value = 15000
$ws.Cells.Item(630, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(630, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(630, 16).Value = 250
$ws.Cells.Item(630, 17).Value = 60
$ws.Cells.Item(630, 18).Value = "Hortaliza"
